$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 40 and 41 swap coin identity (B, C) along with new D/E values.
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"

# E column (Volume 1h %) updates - always safe as text due to padding spaces.
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -1.86%  "
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("E19").Value = "  +9.19%  "
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("E22").Value = "  +19.34%  "
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("E24").Value = "  -7.56%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("E29").Value = "  +4.22%  "
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("E32").Value = "  +7.80%  "
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("E34").Value = "  -0.26%  "
$ws.Range("E35").Value = "  -4.73%  "
$ws.Range("E36").Value = "  -7.33%  "
$ws.Range("E37").Value = "  -5.43%  "
$ws.Range("E38").Value = "  -2.46%  "
$ws.Range("E39").Value = "  -6.35%  "
$ws.Range("E40").Value = "  -3.30%  "
$ws.Range("E41").Value = "  +4.23%  "
$ws.Range("E42").Value = "  -3.06%  "
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("E44").Value = "  +0.52%  "
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("E47").Value = "  +4.62%  "
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("E49").Value = "  -0.99%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  -1.71%  "

# D column (Price) updates.
# Values that are NOT valid numbers (contain multiple dots, e.g. "41.932.12")
# stay text automatically. Values that ARE valid numbers need to be forced to
# text (matching the original inlineStr storage) via a temporary Text format,
# then the style is reset back to Normal so no stray formatting remains.
$ws.Range("D2").Value = "42.008.53"
$ws.Range("D3").Value = "2.210.71"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.28"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.87"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.76"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.06"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Value = "2.533.63"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.16"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.837"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "2.197.60"
$ws.Range("D18").Value = "41.823.73"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000108"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.62"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.11"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.28"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.14"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.49"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.58"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.25"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.03"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.51"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.58"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0789"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.125"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.01"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.26"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.92"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.12"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "64.88"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.60"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.198"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.69"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.70"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.69"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "2.410.61"
